$d = $word.ActiveDocument

# =====================================================================
# Hunk 1: "...seeking an internship with iRobot for the summer of 2018
#          in robotic algorithm development."
#       -> "...seeking an internship with Magna for the summer of 2018
#          in automation."  (with the _GoBack bookmark relocated to sit
#          between "automation" and the final ".")
# =====================================================================

# --- 1a. iRobot -> Magna, fenced so it stays an isolated run ---
$rngI = $d.Content
$rngI.Find.Execute("iRobot") | Out-Null
$iStart = $rngI.Start
$iEnd = $rngI.End
$d.Bookmarks.Add("ZZFenceL", $d.Range($iStart, $iStart)) | Out-Null
$d.Bookmarks.Add("ZZFenceR", $d.Range($iEnd, $iEnd)) | Out-Null

$rngI2 = $d.Content
$rngI2.Find.Execute("iRobot") | Out-Null
$rngI2.Text = "Magna"

$d.Bookmarks.Item("ZZFenceL").Delete()
$d.Bookmarks.Item("ZZFenceR").Delete()

# --- 1b. " robotic algorithm development." -> " automation." ---
# Fence the left edge first so the edit below can't merge backward into
# the preceding " for the summer of 2018 in" run.
$rngDev = $d.Content
$rngDev.Find.Execute(" robotic algorithm development.") | Out-Null
$devStart = $rngDev.Start
$d.Bookmarks.Add("ZZFenceDev", $d.Range($devStart, $devStart)) | Out-Null

$rngDev2 = $d.Content
$rngDev2.Find.Execute(" robotic algorithm development.") | Out-Null
$rngDev2.Text = " automation."

# --- 1c. split " automation." into " " / "automation" / "." ---
$rngCheck = $d.Content
$rngCheck.Find.Execute(" automation.") | Out-Null
$autoStart = $rngCheck.Start
$autoEnd = $rngCheck.End

$splitPos1 = $autoStart + 1        # boundary between " " and "automation"
$splitPos2 = $autoEnd - 1          # boundary between "automation" and "."

$d.Bookmarks.Add("ZZFenceSplit", $d.Range($splitPos1, $splitPos1)) | Out-Null
$d.Bookmarks.Add("ZZFenceDot", $d.Range($splitPos2, $splitPos2)) | Out-Null

# Re-assert the text of the now-isolated "automation" run (change then
# revert) so the engine recomputes xml:space on that run instead of
# blindly inheriting it from the original, space-leading run.
$rngAuto = $d.Range($splitPos1, $splitPos2)
$rngAuto.Text = "automationZZ"
$rngAuto2 = $d.Range($splitPos1, $splitPos1 + 12)
$rngAuto2.Text = "automation"

# Same trick for the trailing "." run.
$dotStart = $splitPos1 + 10
$rngDot = $d.Range($dotStart, $dotStart + 1)
$rngDot.Text = "XX"
$rngDot2 = $d.Range($dotStart, $dotStart + 2)
$rngDot2.Text = "."

# Drop the scaffolding bookmarks (the run splits they created persist).
$d.Bookmarks.Item("ZZFenceDev").Delete()
$d.Bookmarks.Item("ZZFenceSplit").Delete()
$d.Bookmarks.Item("ZZFenceDot").Delete()

# =====================================================================
# Hunk 2: the _GoBack bookmark used to sit between "to track changes
# and fall back on " and "stable versions" -- remove it from there and
# merge the two runs into one.
# =====================================================================

$rngTrack = $d.Content
$rngTrack.Find.Execute("to track changes and fall back on ") | Out-Null
$fenceStart = $rngTrack.Start
$d.Bookmarks.Add("ZZFenceTrack", $d.Range($fenceStart, $fenceStart)) | Out-Null

$d.Bookmarks.Item("_GoBack").Delete()

$rngFull = $d.Content
$rngFull.Find.Execute("to track changes and fall back on stable versions") | Out-Null
$rngFull.Text = "to track changes and fall back on stable versionsZZ"

$rngFull2 = $d.Content
$rngFull2.Find.Execute("to track changes and fall back on stable versionsZZ") | Out-Null
$rngFull2.Text = "to track changes and fall back on stable versions"

$d.Bookmarks.Item("ZZFenceTrack").Delete()

# =====================================================================
# Finally, re-create the _GoBack bookmark at its new home: between
# "automation" and the trailing "." in hunk 1's paragraph.
# =====================================================================
$gbPos = $splitPos1 + 10
$d.Bookmarks.Add("_GoBack", $d.Range($gbPos, $gbPos)) | Out-Null

Write-Output "edit complete"
